# Updates the Price (column D) and Volume(1h) (column E) values for the
# cryptocurrency table in the active worksheet, per the latest data refresh
# (scheduled GitHub Actions run).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# $null for D or E means that column is unchanged for that row and should
# be left untouched.
$updates = @(
    @{Row=2;  D="27.971.28";    E="  -3.27%  "},
    @{Row=3;  D="1.864.39";     E="  -2.25%  "},
    @{Row=4;  D="1.002";        E="  +0.07%  "},
    @{Row=5;  D="318.25";       E="  -1.95%  "},
    @{Row=6;  D="1.001";        E="  +0.01%  "},
    @{Row=7;  D="0.4376";       E="  -4.62%  "},
    @{Row=8;  D="0.3702";       E="  -3.05%  "},
    @{Row=9;  D="0.07511";      E="  -2.62%  "},
    @{Row=10; D="0.9375";       E="  -4.29%  "},
    @{Row=11; D="21.31";        E="  -4.01%  "},
    @{Row=12; D="1.858.14";     E="  -1.86%  "},
    @{Row=13; D="6.737";        E="  -3.15%  "},
    @{Row=14; D="5.451";        E="  -3.99%  "},
    @{Row=15; D="0.06819";      E="  -3.57%  "},
    @{Row=16; D=$null;          E="  +0.00%  "},
    @{Row=17; D="81.64";        E="  -2.61%  "},
    @{Row=18; D="0.000009060";  E="  -4.20%  "},
    @{Row=19; D=$null;          E="  -0.09%  "},
    @{Row=20; D="15.96";        E="  -4.15%  "},
    @{Row=21; D="27.960.52";    E="  -3.30%  "},
    @{Row=22; D=$null;          E="  -3.80%  "},
    @{Row=23; D="11.07";        E="  +1.26%  "},
    @{Row=24; D="2.094.36";     E="  -1.13%  "},
    @{Row=25; D="2.009";        E="  -4.16%  "},
    @{Row=26; D="154.16";       E="  -2.80%  "},
    @{Row=27; D="18.39";        E="  -3.52%  "},
    @{Row=28; D="5.438";        E="  -4.16%  "},
    @{Row=29; D="113.47";       E="  -3.53%  "},
    @{Row=30; D="1.731";        E="  -7.75%  "},
    @{Row=31; D="0.08999";      E="  -3.29%  "},
    @{Row=32; D="0.8119";       E="  -5.90%  "},
    @{Row=33; D="4.822";        E=$null},
    @{Row=34; D="1.175";        E="  -5.60%  "},
    @{Row=35; D="2.923";        E="  -3.44%  "},
    @{Row=36; D=$null;          E="  +0.05%  "},
    @{Row=37; D="0.05500";      E="  -3.71%  "},
    @{Row=38; D="1.119";        E="  -3.37%  "},
    @{Row=39; D="0.01980";      E="  -3.17%  "},
    @{Row=40; D="2.905";        E="  -0.33%  "},
    @{Row=41; D="0.5272";       E="  -4.10%  "},
    @{Row=42; D="7.067";        E="  -5.41%  "},
    @{Row=43; D="0.1693";       E="  -3.49%  "},
    @{Row=44; D="8.812";        E="  -5.71%  "},
    @{Row=45; D="0.06781";      E="  -1.51%  "},
    @{Row=46; D="0.4913";       E="  -5.21%  "},
    @{Row=47; D="10.64";        E="  -5.20%  "},
    @{Row=48; D="106.54";       E="  -3.64%  "},
    @{Row=49; D="1.682";        E="  -5.51%  "},
    @{Row=50; D="0.9997";       E="  -0.06%  "},
    @{Row=51; D="1.901";        E="  -12.50%  "}
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        # The Price column holds plain text (e.g. thousand-grouped values
        # like "27.971.28"). Values that would otherwise parse as an
        # ordinary number (e.g. "1.002", a single decimal point) need an
        # apostrophe prefix so Excel keeps them as text instead of silently
        # converting them to a numeric value. Values with more than one "."
        # (e.g. "27.971.28") are unambiguous and can be written as-is.
        $dVal = $u.D
        if ($dVal -match "^[0-9]+\.[0-9]+$") {
            $ws.Cells.Item($u.Row, 4).Value = "'" + $dVal
        } else {
            $ws.Cells.Item($u.Row, 4).Value = $dVal
        }
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}
